# Update the Prediction (D) and Error (E) columns — and one Cross Entropy
# Loss (F11) value — to the refreshed model-evaluation numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.9999999999999998
$ws.Range("E2").Value = 0.9999999999999998

$ws.Range("D3").Value = 0.9993376189073055
$ws.Range("E3").Value = 0.9993376189073055

$ws.Range("D4").Value = 0.9999999991080499
$ws.Range("E4").Value = 0.9999999991080499

$ws.Range("D5").Value = 0.8671673055247027
$ws.Range("E5").Value = 0.8671673055247027

$ws.Range("D6").Value = 0.9945180849568553
$ws.Range("E6").Value = 0.9945180849568553

$ws.Range("D7").Value = [double]"1.4634860075838E-18"

$ws.Range("D8").Value = 0.9999729626206506
$ws.Range("E8").Value = [double]"2.703737934939276E-05"

$ws.Range("D9").Value = 0.8126159529786328
$ws.Range("E9").Value = 0.1873840470213672

$ws.Range("D10").Value = [double]"4.086794851231971E-16"
$ws.Range("E10").Value = 0.9999999999999996

$ws.Range("D11").Value = [double]"3.617442584425506E-05"
$ws.Range("E11").Value = 0.9999638255741558
$ws.Range("F11").Value = 15.87427997589111
